$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the added "Cap. Percent (%)" column
$ws.Range("G1").Value = "Cap. Percent (%)"

# New formula column: percent of captures = (D2 / $D$3) * 100
$ws.Range("G2").Formula = '=(D2/$D$3)*100'

# New total for column D (abundance*xn sum), mirroring the existing C3 sum
$ws.Range("D3").Formula = "=SUM(D2:D2)"

# Match column widths for the new column G with column F
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth

# Update the last selected cell to reflect the final state recorded in the diff
$ws.Range("F25").Select()

$wb.Save()
